# Deprecate several chart types: re-layout the regenerated grattan chart
# (group "Content Placeholder 3" on slide 1). The chart was rebuilt with a
# slightly narrower plot-area width (cx 10920518 -> 10920451 EMU) and a
# slightly different height (cy 3812069 -> 3812052 EMU), shifting every
# child shape's position/size by a handful of EMU. Apply the new
# Left/Top/Width/Height (points = EMU/12700) to each affected child shape
# inside the group, matched by its original (stable) shape name.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the group shape that holds the rendered chart.
$group = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Content Placeholder 3") {
        $group = $cand
    }
}

$updates = @(
    @{ Name = "rc4"; Left = 62.51748275756836; Top = 147.8544921875; Width = 859.8780517578125; Height = 300.1615905761719 },
    @{ Name = "pl5"; Left = 62.51748275756836; Top = 409.9878845214844; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl6"; Left = 62.51748275756836; Top = 351.9295349121094; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl7"; Left = 62.51748275756836; Top = 293.8711853027344; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl8"; Left = 62.51748275756836; Top = 235.81283569335938; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl9"; Left = 62.51748275756836; Top = 177.75450134277344; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl10"; Left = 99.0044937133789; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl11"; Left = 298.8785095214844; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl12"; Left = 498.7525329589844; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl13"; Left = 698.6265869140625; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl14"; Left = 898.5006713867188; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl15"; Left = 62.51748275756836; Top = 439.01702880859375; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl16"; Left = 62.51748275756836; Top = 380.95867919921875; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl17"; Left = 62.51748275756836; Top = 322.90032958984375; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl18"; Left = 62.51748275756836; Top = 264.8420715332031; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl19"; Left = 62.51748275756836; Top = 206.78370666503906; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl20"; Left = 62.51748275756836; Top = 148.72535705566406; Width = 859.8780517578125; Height = 0.0 },
    @{ Name = "pl21"; Left = 198.94149780273438; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl22"; Left = 398.8155212402344; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl23"; Left = 598.6895751953125; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pl24"; Left = 798.5636596679688; Top = 147.8544921875; Width = 0.0; Height = 300.1615905761719 },
    @{ Name = "pt25"; Left = 320.2590637207031; Top = 308.6843566894531 },
    @{ Name = "pt26"; Left = 371.2269592285156; Top = 308.6843566894531 },
    @{ Name = "pt27"; Left = 260.296875; Top = 287.78338623046875 },
    @{ Name = "pt28"; Left = 439.1841125488281; Top = 304.0397033691406 },
    @{ Name = "pt29"; Left = 484.15576171875; Top = 335.3912048339844 },
    @{ Name = "pt30"; Left = 488.1532287597656; Top = 342.35821533203125 },
    @{ Name = "pt31"; Left = 510.1393737792969; Top = 386.4825439453125 },
    @{ Name = "pt32"; Left = 434.187255859375; Top = 269.20465087890625 },
    @{ Name = "pt33"; Left = 426.1922912597656; Top = 287.78338623046875 },
    @{ Name = "pt34"; Left = 484.15576171875; Top = 329.5853576660156 },
    @{ Name = "pt35"; Left = 484.15576171875; Top = 345.8416748046875 },
    @{ Name = "pt36"; Left = 610.076416015625; Top = 362.0980529785156 },
    @{ Name = "pt37"; Left = 542.1192626953125; Top = 351.6475830078125 },
    @{ Name = "pt38"; Left = 552.1129150390625; Top = 376.0320739746094 },
    @{ Name = "pt39"; Left = 845.9277954101562; Top = 431.7680358886719 },
    @{ Name = "pt40"; Left = 880.7058715820312; Top = 431.7680358886719 },
    @{ Name = "pt41"; Left = 864.9158325195312; Top = 381.837890625 },
    @{ Name = "pt42"; Left = 236.31198120117188; Top = 176.31134033203125 },
    @{ Name = "pt43"; Left = 119.38567352294922; Top = 199.53465270996094 },
    @{ Name = "pt44"; Left = 163.3579559326172; Top = 158.89385986328125 },
    @{ Name = "pt45"; Left = 289.2785949707031; Top = 302.8785095214844 },
    @{ Name = "pt46"; Left = 500.14569091796875; Top = 372.54852294921875 },
    @{ Name = "pt47"; Left = 483.1564025878906; Top = 376.0320739746094 },
    @{ Name = "pt48"; Left = 564.1054077148438; Top = 398.09417724609375 },
    @{ Name = "pt49"; Left = 565.104736328125; Top = 329.5853576660156 },
    @{ Name = "pt50"; Left = 183.34536743164062; Top = 235.53086853027344 },
    @{ Name = "pt51"; Left = 224.3195343017578; Top = 250.62599182128906 },
    @{ Name = "pt52"; Left = 98.99850463867188; Top = 199.53465270996094 },
    @{ Name = "pt53"; Left = 430.1897888183594; Top = 369.0650634765625 },
    @{ Name = "pt54"; Left = 350.24017333984375; Top = 323.779541015625 },
    @{ Name = "pt55"; Left = 510.1393737792969; Top = 378.3543395996094 },
    @{ Name = "pt56"; Left = 352.2389221191406; Top = 304.0397033691406 },
    @{ Name = "tx57"; Left = 47.79393768310547; Top = 435.8634948730469; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "tx58"; Left = 47.79393768310547; Top = 377.8051452636719; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "tx59"; Left = 47.79393768310547; Top = 319.7467956542969; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "tx60"; Left = 47.79393768310547; Top = 261.6884460449219; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "tx61"; Left = 47.79393768310547; Top = 203.6300811767578; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "tx62"; Left = 47.79393768310547; Top = 145.57174682617188; Width = 9.791969299316406; Height = 6.307165622711182 },
    @{ Name = "pl63"; Left = 59.77771759033203; Top = 439.01702880859375 },
    @{ Name = "pl64"; Left = 59.77771759033203; Top = 380.95867919921875 },
    @{ Name = "pl65"; Left = 59.77771759033203; Top = 322.90032958984375 },
    @{ Name = "pl66"; Left = 59.77771759033203; Top = 264.8420715332031 },
    @{ Name = "pl67"; Left = 59.77771759033203; Top = 206.78370666503906 },
    @{ Name = "pl68"; Left = 59.77771759033203; Top = 148.72535705566406 },
    @{ Name = "pl69"; Left = 198.94149780273438; Top = 448.0160827636719 },
    @{ Name = "pl70"; Left = 398.8155212402344; Top = 448.0160827636719 },
    @{ Name = "pl71"; Left = 598.6895751953125; Top = 448.0160827636719 },
    @{ Name = "pl72"; Left = 798.5636596679688; Top = 448.0160827636719 },
    @{ Name = "tx73"; Left = 196.49346923828125; Top = 452.94757080078125; Width = 4.895984649658203; Height = 6.307165622711182 },
    @{ Name = "tx74"; Left = 396.36749267578125; Top = 452.94757080078125; Width = 4.895984649658203; Height = 6.307165622711182 },
    @{ Name = "tx75"; Left = 596.2415771484375; Top = 452.94757080078125; Width = 4.895984649658203; Height = 6.307165622711182 },
    @{ Name = "tx76"; Left = 796.1156005859375; Top = 452.94757080078125; Width = 4.895984649658203; Height = 6.307165622711182 },
    @{ Name = "tx77"; Left = 486.9593200683594; Top = 463.9456787109375; Width = 10.994331359863281; Height = 7.891181468963623 },
    @{ Name = "tx78"; Left = 27.972206115722656; Top = 293.98968505859375; Width = 21.40559196472168; Height = 7.891181468963623 }
)

# Index the updates by shape name for quick lookup.
$byName = @{}
foreach ($u in $updates) {
    $byName[$u.Name] = $u
}

for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $item = $group.GroupItems.Item($i)
    $u = $byName[$item.Name]
    if ($u -ne $null) {
        $item.Left = $u.Left
        $item.Top = $u.Top
        if ($u.ContainsKey("Width")) {
            $item.Width = $u.Width
        }
        if ($u.ContainsKey("Height")) {
            $item.Height = $u.Height
        }
    }
}
